$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")

$ws1.Cells.Item(2, 3).Value = 5
$ws1.Cells.Item(2, 4).Value = 855.12
$ws1.Cells.Item(2, 5).Value = 177.2
$ws1.Cells.Item(3, 3).Value = 5
$ws1.Cells.Item(3, 4).Value = 753.18
$ws1.Cells.Item(3, 5).Value = 152.17
$ws1.Cells.Item(4, 3).Value = 5
$ws1.Cells.Item(4, 4).Value = 729.23
$ws1.Cells.Item(4, 5).Value = 146.01
$ws1.Cells.Item(5, 3).Value = 5
$ws1.Cells.Item(5, 4).Value = 727.32
$ws1.Cells.Item(5, 5).Value = 145.79
$ws1.Cells.Item(6, 3).Value = 5
$ws1.Cells.Item(6, 4).Value = 584.53
$ws1.Cells.Item(6, 5).Value = 118.11
$ws1.Cells.Item(7, 3).Value = 5
$ws1.Cells.Item(7, 4).Value = 579.05
$ws1.Cells.Item(7, 5).Value = 116.37
$ws1.Cells.Item(8, 3).Value = 5
$ws1.Cells.Item(8, 4).Value = 471.26
$ws1.Cells.Item(8, 5).Value = 94.27
$ws1.Cells.Item(15, 3).Value = 1
$ws1.Cells.Item(15, 4).Value = 23.34
$ws1.Cells.Item(15, 5).Value = -5.6
$ws1.Cells.Item(16, 1).Value = 'SUCRIVOIRE (SCRC)'
$ws1.Cells.Item(16, 4).Value = 14.47
$ws1.Cells.Item(16, 5).Value = 7.2
$ws1.Cells.Item(17, 1).Value = 'SICABLE CI (CABC)'
$ws1.Cells.Item(17, 4).Value = 11.95
$ws1.Cells.Item(17, 5).Value = 7.34
$ws1.Cells.Item(18, 1).Value = 'SAFCA CI (SAFC)'
$ws1.Cells.Item(18, 2).Value = 2
$ws1.Cells.Item(18, 4).Value = 11.31
$ws1.Cells.Item(18, 5).Value = 5.57
$ws1.Cells.Item(19, 1).Value = 'UNIWAX CI (UNXC)'
$ws1.Cells.Item(19, 4).Value = 7.26
$ws1.Cells.Item(19, 5).Value = 7.26
$ws1.Cells.Item(20, 1).Value = 'CORIS BANK INTERNATIONAL (CBIBF)'
$ws1.Cells.Item(20, 4).Value = 7.11
$ws1.Cells.Item(20, 5).Value = 7.11
$ws1.Cells.Item(21, 1).Value = 'NSIA BANQUE COTE D''IVOIRE (NSBC)'
$ws1.Cells.Item(21, 4).Value = 5.04
$ws1.Cells.Item(21, 5).Value = 5.04
$ws1.Cells.Item(22, 1).Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Cells.Item(22, 2).Value = 3
$ws1.Cells.Item(22, 3).Value = 2
$ws1.Cells.Item(22, 4).Value = 4.75
$ws1.Cells.Item(22, 5).Value = 4.35
$ws1.Cells.Item(22, 7).Value = '👀 À surveiller'
$ws1.Cells.Item(23, 1).Value = 'SITAB CI (STBC)'
$ws1.Cells.Item(23, 4).Value = 3.68
$ws1.Cells.Item(23, 5).Value = 3.68
$ws1.Cells.Item(24, 1).Value = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$ws1.Cells.Item(24, 3).Value = 0
$ws1.Cells.Item(24, 4).Value = 3.35
$ws1.Cells.Item(24, 5).Value = 3.35
$ws1.Cells.Item(24, 7).Value = '➖ Neutre'
$ws1.Cells.Item(25, 1).Value = 'ONATEL BF (ONTBF)'
$ws1.Cells.Item(25, 3).Value = 1
$ws1.Cells.Item(25, 4).Value = 3.23
$ws1.Cells.Item(25, 5).Value = -1.15
$ws1.Cells.Item(25, 7).Value = '👀 À surveiller'
$ws1.Cells.Item(26, 1).Value = 'CFAO MOTORS CI (CFAC)'
$ws1.Cells.Item(26, 3).Value = 3
$ws1.Cells.Item(26, 4).Value = 2.29
$ws1.Cells.Item(26, 5).Value = 7.19
$ws1.Cells.Item(27, 1).Value = 'VIVO ENERGY CI (SHEC)'
$ws1.Cells.Item(27, 2).Value = 1
$ws1.Cells.Item(27, 3).Value = 1
$ws1.Cells.Item(27, 4).Value = 2.02
$ws1.Cells.Item(27, 5).Value = -2.44
$ws1.Cells.Item(28, 1).Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws1.Cells.Item(28, 3).Value = 2
$ws1.Cells.Item(28, 4).Value = 1.99
$ws1.Cells.Item(28, 5).Value = 7.43
$ws1.Cells.Item(29, 1).Value = 'UNILEVER CI (UNLC)'
$ws1.Cells.Item(29, 2).Value = 1
$ws1.Cells.Item(29, 4).Value = 1.24
$ws1.Cells.Item(29, 5).Value = 7.49
$ws1.Cells.Item(29, 7).Value = '👀 À surveiller'
$ws1.Cells.Item(30, 1).Value = 'SMB CI (SMBC)'
$ws1.Cells.Item(30, 2).Value = 1
$ws1.Cells.Item(30, 4).Value = 0.84
$ws1.Cells.Item(30, 5).Value = -2.15
$ws1.Cells.Item(30, 7).Value = '👀 À surveiller'
$ws1.Cells.Item(31, 1).Value = 'EVIOSYS PACKAGING SIEM CI (SEMC)'
$ws1.Cells.Item(31, 2).Value = 1
$ws1.Cells.Item(31, 4).Value = -0.32
$ws1.Cells.Item(31, 5).Value = -6.81
$ws1.Cells.Item(31, 7).Value = '👀 À surveiller'
$ws1.Cells.Item(32, 1).Value = 'BICI CI (BICC)'
$ws1.Cells.Item(32, 4).Value = -0.49
$ws1.Cells.Item(32, 5).Value = -0.49
$ws1.Cells.Item(33, 1).Value = 'NESTLE CI (NTLC)'
$ws1.Cells.Item(33, 4).Value = -0.89
$ws1.Cells.Item(33, 5).Value = -0.89
$ws1.Cells.Item(34, 1).Value = 'SAPH CI (SPHC)'
$ws1.Cells.Item(34, 4).Value = -0.98
$ws1.Cells.Item(34, 5).Value = -0.98
$ws1.Cells.Item(35, 1).Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Cells.Item(35, 4).Value = -1.14
$ws1.Cells.Item(35, 5).Value = -1.14
$ws1.Cells.Item(36, 1).Value = 'SOGB CI (SOGC)'
$ws1.Cells.Item(36, 4).Value = -2.04
$ws1.Cells.Item(36, 5).Value = -2.04
$ws1.Cells.Item(37, 1).Value = 'ORAGROUP TOGO (ORGT)'
$ws1.Cells.Item(37, 4).Value = -2.08
$ws1.Cells.Item(37, 5).Value = -2.08
$ws1.Cells.Item(38, 1).Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Cells.Item(38, 3).Value = 1
$ws1.Cells.Item(38, 4).Value = -2.89
$ws1.Cells.Item(38, 5).Value = -2.89
$ws1.Cells.Item(39, 1).Value = 'SETAO CI (STAC)'
$ws1.Cells.Item(39, 4).Value = -4
$ws1.Cells.Item(39, 5).Value = -4
$ws1.Cells.Item(40, 1).Value = 'SOLIBRA CI (SLBC)'
$ws1.Cells.Item(40, 3).Value = 1
$ws1.Cells.Item(40, 4).Value = -5.22
$ws1.Cells.Item(40, 5).Value = -5.22
$ws1.Cells.Item(41, 1).Value = 'SICOR CI (SICC)'
$ws1.Cells.Item(41, 2).Value = 0
$ws1.Cells.Item(41, 3).Value = 1
$ws1.Cells.Item(41, 4).Value = -5.26
$ws1.Cells.Item(41, 5).Value = -5.26
$ws1.Cells.Item(41, 6).Value = '🟡 Observer'
$ws1.Cells.Item(41, 7).Value = '➖ Neutre'
$ws1.Cells.Item(42, 1).Value = 'FILTISAC CI (FTSC)'
$ws1.Cells.Item(42, 2).Value = 0
$ws1.Cells.Item(42, 3).Value = 2
$ws1.Cells.Item(42, 4).Value = -6.53
$ws1.Cells.Item(42, 5).Value = -2.86
$ws1.Cells.Item(42, 6).Value = '🟡 Observer'
$ws1.Cells.Item(42, 7).Value = '➖ Neutre'

$ws2 = $wb.Worksheets.Item("Top_YTD")
$ws2.Cells.Item(2, 2).Value = 14517.59
$ws2.Cells.Item(3, 2).Value = 9790.18
$ws2.Cells.Item(4, 2).Value = 8880.79
$ws2.Cells.Item(5, 2).Value = 8811.21
$ws2.Cells.Item(6, 2).Value = 4700.72
$ws2.Cells.Item(7, 2).Value = 4580.92
$ws2.Cells.Item(8, 2).Value = 2665.83
